$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 2867.3635
$ws.Range("I116").Value = 1978.2
$ws.Range("J116").Value = 3608.3333
$ws.Range("K116").Value = 1978.2
$ws.Range("L116").Value = 3608.3333
$ws.Range("M116").Value = 1463.8
$ws.Range("N116").Value = -10492.3333

$ws.Range("H121").Value = 2079.8333
$ws.Range("J121").Value = 2079.8333
$ws.Range("L121").Value = 6239.499899999999
$ws.Range("N121").Value = -9733.499899999999

$ws.Range("H129").Value = 1056.0217
$ws.Range("I129").Value = 1750.4445
$ws.Range("J129").Value = 980.7229
$ws.Range("K129").Value = 5251.333500000001
$ws.Range("L129").Value = 2942.1687
$ws.Range("M129").Value = -251.3335000000006
$ws.Range("N129").Value = -12942.1687

$ws.Range("H135").Value = 25001092
$ws.Range("I135").Value = 1166.5
$ws.Range("J135").Value = 250000420
$ws.Range("K135").Value = 10498.5
$ws.Range("L135").Value = 2250003780
$ws.Range("M135").Value = -7963.5
$ws.Range("N135").Value = -2250008850

$ws.Range("H141").Value = 2573.96
$ws.Range("I141").Value = 1360.8235
$ws.Range("J141").Value = 5151.875
$ws.Range("K141").Value = 4082.4705
$ws.Range("L141").Value = 15455.625
$ws.Range("M141").Value = 1097.5295
$ws.Range("N141").Value = -25815.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3347.1035
$ws.Range("I32").Value = 1474.7778
$ws.Range("J32").Value = 28623.5
$ws.Range("K32").Value = 1474.7778
$ws.Range("L32").Value = 28623.5
$ws.Range("M32").Value = -1187.7778
$ws.Range("N32").Value = -29197.5

$ws.Range("H97").Value = 1042.5217
$ws.Range("I97").Value = 953.0454999999999
$ws.Range("J97").Value = 3011
$ws.Range("K97").Value = 953.0454999999999
$ws.Range("L97").Value = 3011
$ws.Range("M97").Value = -457.0454999999999
$ws.Range("N97").Value = -4003

$ws.Range("H122").Value = 1774.8
$ws.Range("I122").Value = 1418.5
$ws.Range("J122").Value = 3200
$ws.Range("K122").Value = 4255.5
$ws.Range("L122").Value = 9600
$ws.Range("M122").Value = -1805.5
$ws.Range("N122").Value = -14500

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2348.375
$ws.Range("I86").Value = 2296
$ws.Range("J86").Value = 2435.6667
$ws.Range("K86").Value = 2296
$ws.Range("L86").Value = 2435.6667
$ws.Range("M86").Value = -1173
$ws.Range("N86").Value = -4681.6667

$ws.Range("H89").Value = 2348.375
$ws.Range("I89").Value = 2296
$ws.Range("J89").Value = 2435.6667
$ws.Range("K89").Value = 11480
$ws.Range("L89").Value = 12178.3335
$ws.Range("M89").Value = -5864
$ws.Range("N89").Value = -23410.3335

$ws.Range("H96").Value = 16350.385
$ws.Range("I96").Value = 4185
$ws.Range("K96").Value = 4185
$ws.Range("M96").Value = -1439

$ws.Range("H97").Value = 18333.334
$ws.Range("I97").Value = 15000
$ws.Range("K97").Value = 15000
$ws.Range("M97").Value = -14009

$ws.Range("H99").Value = 2512.3076
$ws.Range("I99").Value = 2542.1667
$ws.Range("J99").Value = 2445.125
$ws.Range("K99").Value = 2542.1667
$ws.Range("L99").Value = 2445.125
$ws.Range("M99").Value = -1044.1667
$ws.Range("N99").Value = -5441.125

$ws.Range("H105").Value = 3231.0625
$ws.Range("I105").Value = 2164.5
$ws.Range("J105").Value = 4297.625
$ws.Range("K105").Value = 2164.5
$ws.Range("L105").Value = 4297.625
$ws.Range("M105").Value = -417.5
$ws.Range("N105").Value = -7791.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5152.2812
$ws.Range("I31").Value = 2431.4546
$ws.Range("J31").Value = 5716.981
$ws.Range("K31").Value = 2431.4546
$ws.Range("L31").Value = 5716.981
$ws.Range("M31").Value = -2136.4546
$ws.Range("N31").Value = -6306.981

$ws.Range("H34").Value = 5152.2812
$ws.Range("I34").Value = 2431.4546
$ws.Range("J34").Value = 5716.981
$ws.Range("K34").Value = 2431.4546
$ws.Range("L34").Value = 5716.981
$ws.Range("M34").Value = -2229.4546
$ws.Range("N34").Value = -6120.981

$ws.Range("H103").Value = 15258.909
$ws.Range("I103").Value = 6962
$ws.Range("J103").Value = 20000
$ws.Range("K103").Value = 6962
$ws.Range("L103").Value = 20000
$ws.Range("M103").Value = -5790
$ws.Range("N103").Value = -22344

$ws.Range("H107").Value = 531.3333
$ws.Range("I107").Value = 439.2258
$ws.Range("J107").Value = 790.9091
$ws.Range("K107").Value = 439.2258
$ws.Range("L107").Value = 790.9091
$ws.Range("M107").Value = 1480.7742
$ws.Range("N107").Value = -4630.9091

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 873.54
$ws.Range("I131").Value = 655
$ws.Range("K131").Value = 1965
$ws.Range("M131").Value = 3075

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 800.6667
$ws.Range("I102").Value = 760.8
$ws.Range("K102").Value = 760.8
$ws.Range("M102").Value = 861.2

$ws.Range("H122").Value = 1694
$ws.Range("I122").Value = 1718.625
$ws.Range("J122").Value = 1300
$ws.Range("K122").Value = 5155.875
$ws.Range("L122").Value = 3900
$ws.Range("M122").Value = -2705.875
$ws.Range("N122").Value = -8800

$ws.Range("H132").Value = 32262826
$ws.Range("I132").Value = 62505132
$ws.Range("K132").Value = 187515396
$ws.Range("M132").Value = -187512866

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2307.4092
$ws.Range("I61").Value = 2416.4546
$ws.Range("J61").Value = 2198.3635
$ws.Range("K61").Value = 2416.4546
$ws.Range("L61").Value = 2198.3635
$ws.Range("M61").Value = -2214.4546
$ws.Range("N61").Value = -2602.3635

$ws.Range("H113").Value = 2307.4092
$ws.Range("I113").Value = 2416.4546
$ws.Range("J113").Value = 2198.3635
$ws.Range("K113").Value = 2416.4546
$ws.Range("L113").Value = 2198.3635
$ws.Range("M113").Value = -246.4546
$ws.Range("N113").Value = -6538.363499999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 750.25
$ws.Range("I113").Value = 800.3333
$ws.Range("K113").Value = 2400.9999
$ws.Range("M113").Value = -230.9998999999998

$ws.Range("H126").Value = 1284.1282
$ws.Range("I126").Value = 922.89655
$ws.Range("J126").Value = 2331.7
$ws.Range("K126").Value = 2768.68965
$ws.Range("L126").Value = 6995.099999999999
$ws.Range("M126").Value = -298.6896500000003
$ws.Range("N126").Value = -11935.1
